$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15219.8125
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15936

$ws.Range("H23").Value = 15219.8125
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15468

$ws.Range("H29").Value = 2460.2856
$ws.Range("I29").Value = 5611
$ws.Range("J29").Value = 1200
$ws.Range("K29").Value = 16833
$ws.Range("L29").Value = 3600
$ws.Range("M29").Value = -16552
$ws.Range("N29").Value = -4162

$ws.Range("H38").Value = 868.9167
$ws.Range("J38").Value = 3003
$ws.Range("L38").Value = 9009
$ws.Range("N38").Value = -9753

$ws.Range("H43").Value = 1500
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -2138

$ws.Range("H51").Value = 2984.4443
$ws.Range("I51").Value = 2420
$ws.Range("J51").Value = 3266.6667
$ws.Range("K51").Value = 2420
$ws.Range("L51").Value = 3266.6667
$ws.Range("M51").Value = -1936
$ws.Range("N51").Value = -4234.6667

$ws.Range("H58").Value = 30
$ws.Range("I58").Value = 30
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 90
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 60
$ws.Range("N58").ClearContents()

$ws.Range("H61").Value = 79
$ws.Range("J61").Value = 100
$ws.Range("L61").Value = 300
$ws.Range("N61").Value = -644

$ws.Range("H86").Value = 2199.6667
$ws.Range("I86").Value = 1499
$ws.Range("J86").Value = 2550
$ws.Range("K86").Value = 1499
$ws.Range("L86").Value = 2550
$ws.Range("M86").Value = -376
$ws.Range("N86").Value = -4796

$ws.Range("H87").Value = 25395
$ws.Range("J87").Value = 25966.666
$ws.Range("L87").Value = 25966.666
$ws.Range("N87").Value = -28462.666

$ws.Range("H89").Value = 2199.6667
$ws.Range("I89").Value = 1499
$ws.Range("J89").Value = 2550
$ws.Range("K89").Value = 7495
$ws.Range("L89").Value = 12750
$ws.Range("M89").Value = -1879
$ws.Range("N89").Value = -23982

$ws.Range("H90").Value = 25395
$ws.Range("J90").Value = 25966.666
$ws.Range("L90").Value = 77899.99800000001
$ws.Range("N90").Value = -90379.99800000001

$ws.Range("H112").Value = 6994104
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 7793361.5
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 23380084.5
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -23382300.5

$ws.Range("H137").Value = 1281.125
$ws.Range("J137").Value = 1616.5518
$ws.Range("L137").Value = 4849.6554
$ws.Range("N137").Value = -9949.6554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3345
$ws.Range("I25").Value = 3345
$ws.Range("K25").Value = 3345
$ws.Range("M25").Value = -2943

$ws.Range("H31").Value = 13320
$ws.Range("I31").Value = 2866.6667
$ws.Range("J31").Value = 29000
$ws.Range("K31").Value = 2866.6667
$ws.Range("L31").Value = 29000
$ws.Range("M31").Value = -2572.6667
$ws.Range("N31").Value = -29588

$ws.Range("H35").Value = 887
$ws.Range("I35").Value = 887
$ws.Range("K35").Value = 887
$ws.Range("M35").Value = -481

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1056.8334
$ws.Range("I37").Value = 1028.2
$ws.Range("J37").Value = 1200
$ws.Range("K37").Value = 1028.2
$ws.Range("L37").Value = 1200
$ws.Range("M37").Value = -891.2
$ws.Range("N37").Value = -1474

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 54166.332
$ws.Range("J20").Value = 54166.332
$ws.Range("L20").Value = 54166.332
$ws.Range("N20").Value = -54638.332

$ws.Range("H30").Value = 54166.332
$ws.Range("J30").Value = 54166.332
$ws.Range("L30").Value = 54166.332
$ws.Range("N30").Value = -54348.332

$ws.Range("H128").Value = 54166.332
$ws.Range("J128").Value = 54166.332
$ws.Range("L128").Value = 54166.332
$ws.Range("N128").Value = -64126.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3041.3833
$ws.Range("I68").Value = 3160.4424
$ws.Range("J68").Value = 2267.5
$ws.Range("K68").Value = 9481.3272
$ws.Range("L68").Value = 6802.5
$ws.Range("M68").Value = -8670.3272
$ws.Range("N68").Value = -8424.5

$ws.Range("H71").Value = 3041.3833
$ws.Range("I71").Value = 3160.4424
$ws.Range("J71").Value = 2267.5
$ws.Range("K71").Value = 28443.9816
$ws.Range("L71").Value = 20407.5
$ws.Range("M71").Value = -24387.9816
$ws.Range("N71").Value = -28519.5

$ws.Range("H92").Value = 696.6667
$ws.Range("I92").Value = 696.6667
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2090.0001
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -842.0001000000002
$ws.Range("N92").ClearContents()

$ws.Range("H131").Value = 2370.4482
$ws.Range("J131").Value = 2614.039
$ws.Range("L131").Value = 7842.117
$ws.Range("N131").Value = -17922.117

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2481
$ws.Range("I80").Value = 2401.25
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 2401.25
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -1403.25
$ws.Range("N80").Value = -4796

$ws.Range("H83").Value = 2481
$ws.Range("I83").Value = 2401.25
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 12006.25
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -7014.25
$ws.Range("N83").Value = -23984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 13135.333
$ws.Range("J74").Value = 13248.6
$ws.Range("L74").Value = 13248.6
$ws.Range("N74").Value = -15120.6

$ws.Range("H77").Value = 13135.333
$ws.Range("J77").Value = 13248.6
$ws.Range("L77").Value = 39745.8
$ws.Range("N77").Value = -49105.8

$ws.Range("H136").Value = 14494193
$ws.Range("I136").Value = 20834018
$ws.Range("K136").Value = 62502054
$ws.Range("M136").Value = -62499504
